$d = $word.ActiveDocument

# Work without tracked-changes markup being injected into the runs we touch.
$origTrackRevisions = $d.TrackRevisions
$d.TrackRevisions = $false

# ---------------------------------------------------------------------------
# 1) Big centered title ("...BREATING CONSIDERATIONS PREVENTION" @ sz=52):
#    fix the "BREATING" -> "BREATHING" typo and re-split the run that used
#    to hold "CONSIDERATIONS" into its own trailing-space run, matching the
#    target run layout exactly (including the rsid carried on the 2nd/3rd
#    runs, which come from the original "CONSIDERATIONS"/" " runs).
# ---------------------------------------------------------------------------
$titlePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "BREATING CONSIDERATIONS PREVENTION*") {
        $titlePara = $p
        break
    }
}

if ($titlePara -ne $null) {
    $titleXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="0ED4CEE7" w14:textId="0EB86F7D" w:rsidR="00446DFC" w:rsidRDefault="009B58CF" w:rsidP="00B111EA"><w:pPr><w:jc w:val="center"/><w:rPr><w:bCs/><w:sz w:val="52"/><w:szCs w:val="44"/></w:rPr></w:pPr><w:r><w:rPr><w:bCs/><w:sz w:val="52"/><w:szCs w:val="44"/></w:rPr><w:t>BREAT</w:t></w:r><w:r w:rsidR="0080495A"><w:rPr><w:bCs/><w:sz w:val="52"/><w:szCs w:val="44"/></w:rPr><w:t>H</w:t></w:r><w:r w:rsidR="0080495A"><w:rPr><w:bCs/><w:sz w:val="52"/><w:szCs w:val="44"/></w:rPr><w:t xml:space="preserve">ING </w:t></w:r><w:r><w:rPr><w:bCs/><w:sz w:val="52"/><w:szCs w:val="44"/></w:rPr><w:t xml:space="preserve">CONSIDERATIONS </w:t></w:r><w:r w:rsidR="00446DFC"><w:rPr><w:bCs/><w:sz w:val="52"/><w:szCs w:val="44"/></w:rPr><w:t>PREVENTION</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
    [void]$titlePara.Range.InsertXML($titleXml)
}

# ---------------------------------------------------------------------------
# 2) Bold sub-heading ("BREATHING CONSIDERATIONS PREVENTION SECURITY
#    SYSTEMS" @ sz=24): merge the "CONSIDERATIONS" run and the following
#    " " run into a single "CONSIDERATIONS " run.
# ---------------------------------------------------------------------------
$subheadPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "BREATHING CONSIDERATIONS PREVENTION SECURITY SYSTEMS*") {
        $subheadPara = $p
        break
    }
}

if ($subheadPara -ne $null) {
    $subheadXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="2977E484" w14:textId="1D1232CB" w:rsidR="008D77DA" w:rsidRPr="00C0532F" w:rsidRDefault="009B58CF" w:rsidP="008D77DA"><w:pPr><w:ind w:left="360" w:hanging="360"/><w:jc w:val="both"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="24"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">BREATHING </w:t></w:r><w:r w:rsidR="0080495A"><w:rPr><w:b/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">CONSIDERATIONS </w:t></w:r><w:r w:rsidR="008D77DA"><w:rPr><w:b/><w:sz w:val="24"/></w:rPr><w:t>PREVENTION SECURITY SYSTEMS</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
    [void]$subheadPara.Range.InsertXML($subheadXml)
}

# Restore the document's original track-changes setting.
$d.TrackRevisions = $origTrackRevisions
